$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 630
$ws1.Range("F8").Value = 327
$ws1.Range("F9").Value = 1713
$ws1.Range("F10").Value = 358
$ws1.Range("F11").Value = 1414
$ws1.Range("F12").Value = 805
$ws1.Range("F13").Value = 332
$ws1.Range("F15").Value = 12732
$ws1.Range("F16").Value = 12756
$ws1.Range("F18").Value = 743
$ws1.Range("F19").Value = 10
$ws1.Range("F20").Value = 511
$ws1.Range("F22").Value = 544
$ws1.Range("F25").Value = 15
$ws1.Range("F26").Value = 240

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 78

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 165

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 165
$ws4.Range("F11").Value = 630
$ws4.Range("F13").Value = 327
$ws4.Range("F14").Value = 1713
$ws4.Range("F15").Value = 358
$ws4.Range("F16").Value = 1414
$ws4.Range("F17").Value = 805
$ws4.Range("F18").Value = 332
$ws4.Range("F19").Value = 78
$ws4.Range("F21").Value = 12732
$ws4.Range("F22").Value = 12756
$ws4.Range("F24").Value = 743
$ws4.Range("F25").Value = 10
$ws4.Range("F26").Value = 511
$ws4.Range("F28").Value = 544
$ws4.Range("F34").Value = 15
$ws4.Range("F36").Value = 240
